$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.807.56"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "3.505.89"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.81"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.53"
$ws.Range("E6").Value = "  +4.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +1.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.138"
$ws.Range("E9").Value = "  +3.43%  "

$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "4.110.51"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.38"
$ws.Range("E13").Value = "  +12.75%  "

$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "67.763.70"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("E16").Value = "  +0.81%  "

$ws.Range("D17").Value = "3.504.37"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.74"
$ws.Range("E19").Value = "  +3.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.82"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.44"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("E26").Value = "  +2.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.51"
$ws.Range("E27").Value = "  +2.98%  "

$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.29"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.47"
$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "24.07"
$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.45"
$ws.Range("E34").Value = "  +0.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.68"
$ws.Range("E36").Value = "  +2.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.44"
$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.877"
$ws.Range("E39").Value = "  -1.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.14"
$ws.Range("E40").Value = "  +3.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.74"
$ws.Range("E41").Value = "  +1.22%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  +3.45%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.80"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.65"
$ws.Range("E44").Value = "  +1.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0740"
$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("D46").Value = "2.832.85"
$ws.Range("E46").Value = "  +1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.41"
$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0306"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "345.41"
$ws.Range("E49").Value = "  +1.26%  "

$ws.Range("E50").Value = "  -0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.66"
$ws.Range("E51").Value = "  +0.91%  "
